# Insert one new week of "Cebollín" price data (Primera + Segunda) at the
# top of the existing weekly block (rows 385-405), pushing the rest of the
# block down by 2 rows (385-405 -> 387-407).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 385, shifting everything
# below (including the former row 385) down by two rows.
$ws.Rows.Item(385).Resize(2).Insert()

# --- New row 385: "Primera" quality, week of 2022-01-24 ---
$ws.Cells.Item(385,1).Value  = 3
$ws.Cells.Item(385,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(385,3).Value  = "Coquimbo"
$ws.Cells.Item(385,4).Value  = 44585
$ws.Cells.Item(385,5).Value  = 5
$ws.Cells.Item(385,6).Value  = 100112037
$ws.Cells.Item(385,7).Value  = "Cebollín"
$ws.Cells.Item(385,8).Value  = "Sin especificar"
$ws.Cells.Item(385,9).Value  = "Primera"
$ws.Cells.Item(385,10).Value = 210
$ws.Cells.Item(385,11).Value = 3000
$ws.Cells.Item(385,12).Value = 3500
$ws.Cells.Item(385,13).Value = 3238
$ws.Cells.Item(385,14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(385,15).Value = "Provincia de Quillota"
$ws.Cells.Item(385,16).Value = 90
$ws.Cells.Item(385,17).Value = 36
$ws.Cells.Item(385,18).Value = "Hortaliza"

# --- New row 386: "Segunda" quality, same week ---
$ws.Cells.Item(386,1).Value  = 3
$ws.Cells.Item(386,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(386,3).Value  = "Coquimbo"
$ws.Cells.Item(386,4).Value  = 44585
$ws.Cells.Item(386,5).Value  = 5
$ws.Cells.Item(386,6).Value  = 100112037
$ws.Cells.Item(386,7).Value  = "Cebollín"
$ws.Cells.Item(386,8).Value  = "Sin especificar"
$ws.Cells.Item(386,9).Value  = "Segunda"
$ws.Cells.Item(386,10).Value = 80
$ws.Cells.Item(386,11).Value = 2500
$ws.Cells.Item(386,12).Value = 2500
$ws.Cells.Item(386,13).Value = 2500
$ws.Cells.Item(386,14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(386,15).Value = "Provincia de Quillota"
$ws.Cells.Item(386,16).Value = 69
$ws.Cells.Item(386,17).Value = 36
$ws.Cells.Item(386,18).Value = "Hortaliza"
